$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Remove the stray empty B5 cell in the "ODI Batting" sheet
#    (diff shows the <c r="B5" t="inlineStr"/> element being dropped)
# ------------------------------------------------------------------
$wsOdiBatting = $wb.Worksheets.Item("ODI Batting")
$wsOdiBatting.Cells.Item(5, 2).ClearContents()

# ------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" worksheet after "ODI Batting"
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsExtra.Name = "ODI Batting Extra"
$ws = $wsExtra

# ------------------------------------------------------------------
# 3. Header row
# ------------------------------------------------------------------
$headers = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le 6; $c++) {
    $ws.Cells.Item(1, $c).NumberFormat = "@"
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# ------------------------------------------------------------------
# 4. Data rows
# ------------------------------------------------------------------
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = '4480'
$ws.Cells.Item(2,2).Value = 3
$ws.Cells.Item(2,3).NumberFormat = "@"
$ws.Cells.Item(2,3).Value = '8'
$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = '2'
$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value = '22.43%'
$ws.Cells.Item(2,6).NumberFormat = "@"
$ws.Cells.Item(2,6).Value = 'NO'
$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = '4482'
$ws.Cells.Item(3,2).Value = 3
$ws.Cells.Item(3,3).NumberFormat = "@"
$ws.Cells.Item(3,3).Value = '0'
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = '0'
$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value = '0.36%'
$ws.Cells.Item(3,6).NumberFormat = "@"
$ws.Cells.Item(3,6).Value = 'NO'
$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = '4533'
$ws.Cells.Item(4,2).NumberFormat = "@"
$ws.Cells.Item(4,2).Value = ''
$ws.Cells.Item(4,3).NumberFormat = "@"
$ws.Cells.Item(4,3).Value = ''
$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value = ''
$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value = ''
$ws.Cells.Item(4,6).NumberFormat = "@"
$ws.Cells.Item(4,6).Value = 'NO'
$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = '4637'
$ws.Cells.Item(5,2).Value = 3
$ws.Cells.Item(5,3).NumberFormat = "@"
$ws.Cells.Item(5,3).Value = ''
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = ''
$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value = ''
$ws.Cells.Item(5,6).NumberFormat = "@"
$ws.Cells.Item(5,6).Value = 'NO'
$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = '4640'
$ws.Cells.Item(6,2).Value = 4
$ws.Cells.Item(6,3).NumberFormat = "@"
$ws.Cells.Item(6,3).Value = '0'
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = '0'
$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value = '3.59%'
$ws.Cells.Item(6,6).NumberFormat = "@"
$ws.Cells.Item(6,6).Value = 'NO'
$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = '4643'
$ws.Cells.Item(7,2).NumberFormat = "@"
$ws.Cells.Item(7,2).Value = ''
$ws.Cells.Item(7,3).NumberFormat = "@"
$ws.Cells.Item(7,3).Value = ''
$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value = ''
$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value = ''
$ws.Cells.Item(7,6).NumberFormat = "@"
$ws.Cells.Item(7,6).Value = 'NO'
$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = '4656'
$ws.Cells.Item(8,2).NumberFormat = "@"
$ws.Cells.Item(8,2).Value = ''
$ws.Cells.Item(8,3).NumberFormat = "@"
$ws.Cells.Item(8,3).Value = ''
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = ''
$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value = ''
$ws.Cells.Item(8,6).NumberFormat = "@"
$ws.Cells.Item(8,6).Value = 'NO'
$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = '4657'
$ws.Cells.Item(9,2).Value = 3
$ws.Cells.Item(9,3).NumberFormat = "@"
$ws.Cells.Item(9,3).Value = '4'
$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value = '7'
$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value = '32.98%'
$ws.Cells.Item(9,6).NumberFormat = "@"
$ws.Cells.Item(9,6).Value = 'NO'
$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = '4658'
$ws.Cells.Item(10,2).Value = 3
$ws.Cells.Item(10,3).NumberFormat = "@"
$ws.Cells.Item(10,3).Value = '2'
$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value = '0'
$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value = '9.52%'
$ws.Cells.Item(10,6).NumberFormat = "@"
$ws.Cells.Item(10,6).Value = 'NO'
$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = '4685'
$ws.Cells.Item(11,2).NumberFormat = "@"
$ws.Cells.Item(11,2).Value = ''
$ws.Cells.Item(11,3).NumberFormat = "@"
$ws.Cells.Item(11,3).Value = ''
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = ''
$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value = ''
$ws.Cells.Item(11,6).NumberFormat = "@"
$ws.Cells.Item(11,6).Value = 'NO'
$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = '4692'
$ws.Cells.Item(12,2).NumberFormat = "@"
$ws.Cells.Item(12,2).Value = ''
$ws.Cells.Item(12,3).NumberFormat = "@"
$ws.Cells.Item(12,3).Value = ''
$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value = ''
$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value = ''
$ws.Cells.Item(12,6).NumberFormat = "@"
$ws.Cells.Item(12,6).Value = 'NO'
$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = '4695'
$ws.Cells.Item(13,2).Value = 4
$ws.Cells.Item(13,3).NumberFormat = "@"
$ws.Cells.Item(13,3).Value = '2'
$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value = '0'
$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value = '7.21%'
$ws.Cells.Item(13,6).NumberFormat = "@"
$ws.Cells.Item(13,6).Value = 'NO'
$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = '4697'
$ws.Cells.Item(14,2).Value = 4
$ws.Cells.Item(14,3).NumberFormat = "@"
$ws.Cells.Item(14,3).Value = '1'
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = '1'
$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value = '4.42%'
$ws.Cells.Item(14,6).NumberFormat = "@"
$ws.Cells.Item(14,6).Value = 'NO'
$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = '4725'
$ws.Cells.Item(15,2).Value = 1
$ws.Cells.Item(15,3).NumberFormat = "@"
$ws.Cells.Item(15,3).Value = '0'
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = '0'
$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value = '1.57%'
$ws.Cells.Item(15,6).NumberFormat = "@"
$ws.Cells.Item(15,6).Value = 'NO'

# ------------------------------------------------------------------
# 5. Copy the header formatting (bold font + border + alignment)
#    from "ODI Batting"!A1:F1 so the new header matches the style
#    used by the other sheets in the workbook.
# ------------------------------------------------------------------
$srcHeader = $wsOdiBatting.Range("A1:F1")
$srcHeader.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
